$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell $ws "D2" "318.59"
Set-TextCell $ws "E2" "4.31%"
Set-TextCell $ws "G2" "23"
Set-TextCell $ws "D3" "39.75"
Set-TextCell $ws "E3" "1.75%"
Set-TextCell $ws "G3" "23"
Set-TextCell $ws "D4" "5.147"
Set-TextCell $ws "E4" "0.66%"
Set-TextCell $ws "G4" "23"
Set-TextCell $ws "D5" "0.08229"
Set-TextCell $ws "E5" "1.99%"
Set-TextCell $ws "G5" "23"
Set-TextCell $ws "D6" "2.055"
Set-TextCell $ws "E6" "6.18%"
Set-TextCell $ws "G6" "23"
Set-TextCell $ws "D7" "8.378"
Set-TextCell $ws "E7" "4.12%"
Set-TextCell $ws "G7" "23"
Set-TextCell $ws "D8" "4.319"
Set-TextCell $ws "E8" "2.47%"
Set-TextCell $ws "G8" "23"
Set-TextCell $ws "D9" "0.9400"
Set-TextCell $ws "E9" "1.48%"
Set-TextCell $ws "G9" "23"
Set-TextCell $ws "D10" "0.1356"
Set-TextCell $ws "E10" "-0.25%"
Set-TextCell $ws "G10" "23"
Set-TextCell $ws "D11" "0.1995"
Set-TextCell $ws "E11" "4.24%"
Set-TextCell $ws "G11" "23"
Set-TextCell $ws "D12" "0.09089"
Set-TextCell $ws "E12" "0.99%"
Set-TextCell $ws "G12" "23"
Set-TextCell $ws "D13" "0.03525"
Set-TextCell $ws "E13" "0.16%"
Set-TextCell $ws "G13" "23"
Set-TextCell $ws "D14" "0.09799"
Set-TextCell $ws "E14" "0.40%"
Set-TextCell $ws "G14" "23"
Set-TextCell $ws "D15" "0.001405"
Set-TextCell $ws "E15" "-0.06%"
Set-TextCell $ws "G15" "23"
Set-TextCell $ws "D16" "0.006323"
Set-TextCell $ws "E16" "7.01%"
Set-TextCell $ws "G16" "23"
Set-TextCell $ws "E17" "-1.97%"
Set-TextCell $ws "G17" "23"
Set-TextCell $ws "D18" "3.242"
Set-TextCell $ws "E18" "-3.89%"
Set-TextCell $ws "G18" "23"
Set-TextCell $ws "D19" "0.3496"
Set-TextCell $ws "E19" "1.00%"
Set-TextCell $ws "G19" "23"
Set-TextCell $ws "D20" "0.1322"
Set-TextCell $ws "E20" "0.27%"
Set-TextCell $ws "G20" "23"
Set-TextCell $ws "D21" "4.956"
Set-TextCell $ws "E21" "5.86%"
Set-TextCell $ws "G21" "23"
Set-TextCell $ws "D22" "0.2451"
Set-TextCell $ws "E22" "1.38%"
Set-TextCell $ws "G22" "23"
Set-TextCell $ws "D23" "0.04342"
Set-TextCell $ws "E23" "-0.58%"
Set-TextCell $ws "G23" "23"
Set-TextCell $ws "D24" "0.001234"
Set-TextCell $ws "E24" "2.17%"
Set-TextCell $ws "G24" "23"
Set-TextCell $ws "D25" "0.004797"
Set-TextCell $ws "E25" "12.40%"
Set-TextCell $ws "G25" "23"
Set-TextCell $ws "E26" "-0.07%"
Set-TextCell $ws "G26" "23"
Set-TextCell $ws "D27" "0.0003998"
Set-TextCell $ws "E27" "-10.10%"
Set-TextCell $ws "G27" "23"
Set-TextCell $ws "G28" "23"
Set-TextCell $ws "G29" "23"
Set-TextCell $ws "G30" "23"
Set-TextCell $ws "G31" "23"
Set-TextCell $ws "G32" "23"
Set-TextCell $ws "G33" "23"
Set-TextCell $ws "G34" "23"
Set-TextCell $ws "G35" "23"
Set-TextCell $ws "G36" "23"
Set-TextCell $ws "G37" "23"
Set-TextCell $ws "G38" "23"
Set-TextCell $ws "D39" "0.02328"
Set-TextCell $ws "E39" "14.83%"
Set-TextCell $ws "G39" "23"
Set-TextCell $ws "D40" "0.05209"
Set-TextCell $ws "E40" "3.58%"
Set-TextCell $ws "G40" "23"
Set-TextCell $ws "D41" "0.007760"
Set-TextCell $ws "E41" "3.16%"
Set-TextCell $ws "G41" "23"
Set-TextCell $ws "D42" "0.01014"
Set-TextCell $ws "E42" "4.95%"
Set-TextCell $ws "G42" "23"
Set-TextCell $ws "E43" "5.04%"
Set-TextCell $ws "G43" "23"
Set-TextCell $ws "D44" "0.002042"
Set-TextCell $ws "E44" "-2.52%"
Set-TextCell $ws "G44" "23"
Set-TextCell $ws "D45" "0.009328"
Set-TextCell $ws "E45" "-4.69%"
Set-TextCell $ws "G45" "23"
Set-TextCell $ws "D46" "0.00006612"
Set-TextCell $ws "E46" "6.19%"
Set-TextCell $ws "G46" "23"
Set-TextCell $ws "D47" "0.00000000750"
Set-TextCell $ws "E47" "-0.25%"
Set-TextCell $ws "G47" "23"
Set-TextCell $ws "D48" "0.002945"
Set-TextCell $ws "E48" "2.48%"
Set-TextCell $ws "G48" "23"
Set-TextCell $ws "D49" "0.001689"
Set-TextCell $ws "E49" "-6.33%"
Set-TextCell $ws "G49" "23"
Set-TextCell $ws "D50" "0.00002099"
Set-TextCell $ws "E50" "-0.25%"
Set-TextCell $ws "G50" "23"
Set-TextCell $ws "D51" "0.0001999"
Set-TextCell $ws "E51" "-0.25%"
Set-TextCell $ws "G51" "23"

Write-Host "Applied all crypto price/volume/hour updates"